# pentest, linux cpu info, andriod studio
#
# Appends one new row (29) to the bottom of the TOPIC/DSCP/CODE table on
# sheet 1: "linux" / "check cpu" / "cat /proc/cpuinfo" -- formatted the
# same way as the row above it (wrap-text, Arial 10) -- then moves the
# viewport/selection down to show it, mirroring where Excel's UI would
# land right after typing the new row and pressing Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row, directly below the existing last row (28).
$ws.Range("A29").Value = "linux"
$ws.Range("B29").Value = "check cpu"
$ws.Range("C29").Value = "cat /proc/cpuinfo"

# Match the formatting of the row above (wrap text, same font/style) by
# copying its formats onto the new row, same as Excel does when you
# continue a formatted table downward.
$ws.Range("A28:C28").Copy() | Out-Null
$ws.Range("A29:C29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view: scroll so row 22 is at the top and land the selection
# one row below the new data, same as the post-edit cursor position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A30").Select()
